# Update the LR-pair TPM-derived metrics for the Wnt10b-Fzd7 sheet.
# This mirrors a re-run of the NATMI scoring scripts with new TPM values:
# the receptor average/total expression (M2/N2) changed, which in turn
# changes every downstream "edge weight" / "derived specificity" column
# (O, P, Q, R, S, T) for all three rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Wnt10b -> Fzd7 -> ECs)
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 3.41460036991
$ws.Range("R2").Value = 30.73140332919
$ws.Range("S2").Value = 0.04063212692754557
$ws.Range("T2").Value = 0.04063212692754556

# Row 3 (FAPs -> Wnt10b -> Fzd7 -> FAPs)
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 34.16927966704566
$ws.Range("S3").Value = 0.4065982422683317
$ws.Range("T3").Value = 0.4065982422683317

# Row 4 (FAPs -> Wnt10b -> Fzd7 -> MuSCs)
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 46.45307859922034
$ws.Range("R4").Value = 418.077707392983
$ws.Range("S4").Value = 0.5527696308041227
$ws.Range("T4").Value = 0.5527696308041226
